$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.237273216247559
$ws.Range("B1").Value = 1.431785464286804
$ws.Range("C1").Value = 1.811576008796692
$ws.Range("D1").Value = 3.315759181976318
$ws.Range("E1").Value = 15
